# Updates the cryptos list (Price and Volume(1h) columns, plus a 3-way
# reorder of the NEARProtocol / RenderToken / VeChain rows) to reflect the
# latest scrape, per the "Updated cryptos list" GitHub Actions commit.
#
# Column D (Price) is written as text (NumberFormat "@") because several
# values look numeric (e.g. "1.00", "0.0310", "13.00") and Excel would
# otherwise silently coerce them to numbers and drop the significant
# trailing/leading zeros. The Style is reset to "Normal" right afterwards
# so the cell keeps its original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.842.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -6.74%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.537.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.573'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.89%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.550'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.81'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.31%  '

$ws.Range("E12").Value = '  -4.13%  '

$ws.Range("E13").Value = '  +5.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.923.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.558.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.869'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.864.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.00%  '

$ws.Range("E21").Value = '  -2.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.03%  '

$ws.Range("E24").Value = '  -3.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '29.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.13%  '

$ws.Range("E26").Value = '  -6.78%  '

$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.08%  '

$ws.Range("E30").Value = '  -3.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.81%  '

$ws.Range("E33").Value = '  -5.39%  '

$ws.Range("E34").Value = '  -1.90%  '

$ws.Range("E35").Value = '  -6.07%  '

$ws.Range("E36").Value = '  -5.05%  '

$ws.Range("E37").Value = '  -6.90%  '

$ws.Range("E38").Value = '  -3.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.44%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.22%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0310'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.24%  '

$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.089.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '84.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.781.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.39%  '

$ws.Range("E51").Value = '  -4.54%  '
